$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "exclusivement,TPE,PME,met,place"
$ws.Range("C2").Value = "met en place"
$ws.Range("D2").Value = "un patriotisme économique"
$ws.Range("E2").Value = "un protectionnisme intelligent"
$ws.Range("F2").Value = "protectionnisme,intelligent,dit,constructeurs"

$ws.Range("B4").Value = "mettre,œuvre,protectionnisme,intelligent,mettre"
$ws.Range("C4").Value = "en avant de"
$ws.Range("D4").Value = "le patriotisme économique"
$ws.Range("E4").Value = "pour donner un"
$ws.Range("F4").Value = "donner,avantage,entreprises,françaises"

$ws.Range("B5").Value = "PME,met,place,patriotisme,économique"
$ws.Range("C5").Value = "un patriotisme économique"
$ws.Range("D5").Value = "un protectionnisme intelligent"
$ws.Range("E5").Value = "il dit à"
$ws.Range("F5").Value = "dit,constructeurs,américains,voulez"

$ws.Range("B12").Value = "patriotisme,économique,protectionnisme,intelligent,dit"
$ws.Range("C12").Value = "il dit à"
$ws.Range("D12").Value = "les constructeurs américains"
$ws.Range("E12").Value = "si vous voulez"
$ws.Range("F12").Value = "voulez,aller,faire,voitures"

$ws.Range("B13").Value = "mettre,patriotisme,économique,donner,avantage"
$ws.Range("C13").Value = "un avantage à"
$ws.Range("D13").Value = "les entreprises françaises"
$ws.Range("E13").Value = "dans la commande"
$ws.Range("F13").Value = "commande,publique,patriotisme,économique"

$ws.Range("B14").Value = "supplémentaires,suppression,travail,détaché,baisse"
$ws.Range("C14").Value = "la baisse de"
$ws.Range("D14").Value = "les charges"
$ws.Range("E14").Value = "mais exclusivement pour"
$ws.Range("F14").Value = "exclusivement,TPE,PME,met"

$ws.Range("B17").Value = "charges,exclusivement,TPE,PME,met"
$ws.Range("C17").Value = "Il met en"
$ws.Range("E17").Value = "un patriotisme économique"
$ws.Range("F17").Value = "patriotisme,économique,protectionnisme,intelligent"

$ws.Range("B18").Value = "Trump,intéresse,cause,puisqu,met"
$ws.Range("C18").Value = "il met en"
$ws.Range("E18").Value = "la politique que"
$ws.Range("F18").Value = "politique,appelle,vœux,notamment"
